$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated detection-result data (45 rows x 7 cols: A:G) reflecting the
# revised "kombinasi loc" (bounding-box combination) parameters.
$data = New-Object 'object[,]' 45,7

$data[0,0] = "AK_1.png"; $data[0,1] = 0.394; $data[0,2] = 0.001; $data[0,3] = 0.095; $data[0,4] = 0.367; $data[0,5] = "Tidak Diketahui"; $data[0,6] = "Salah"
$data[1,0] = "AK_2.png"; $data[1,1] = 0.755; $data[1,2] = 0.002; $data[1,3] = 0.337; $data[1,4] = 0.7; $data[1,5] = "Akhlak Kamiswara"; $data[1,6] = "Benar"
$data[2,0] = "AK_3.png"; $data[2,1] = 0.663; $data[2,2] = 0.002; $data[2,3] = 0.427; $data[2,4] = 0.733; $data[2,5] = "Akhlak Kamiswara"; $data[2,6] = "Benar"
$data[3,0] = "AK_4.png"; $data[3,1] = 0.324; $data[3,2] = 0.001; $data[3,3] = 0.43; $data[3,4] = 0.267; $data[3,5] = "Tidak Diketahui"; $data[3,6] = "Salah"
$data[4,0] = "AK_5.png"; $data[4,1] = 0.339; $data[4,2] = 0.001; $data[4,3] = 0.439; $data[4,4] = 0.333; $data[4,5] = "Tidak Diketahui"; $data[4,6] = "Salah"
$data[5,0] = "MIB_1.png"; $data[5,1] = 1.356; $data[5,2] = 0.004; $data[5,3] = 0.08699999999999999; $data[5,4] = 0.5; $data[5,5] = "Muhammad Iqbal Baqi"; $data[5,6] = "Benar"
$data[6,0] = "MIB_2.png"; $data[6,1] = 1.038; $data[6,2] = 0.003; $data[6,3] = 0.275; $data[6,4] = 0.633; $data[6,5] = "Muhammad Iqbal Baqi"; $data[6,6] = "Benar"
$data[7,0] = "MIB_3.png"; $data[7,1] = 1.56; $data[7,2] = 0.005; $data[7,3] = 0.545; $data[7,4] = 0.9; $data[7,5] = "Muhammad Iqbal Baqi"; $data[7,6] = "Benar"
$data[8,0] = "MIB_4.png"; $data[8,1] = 0.983; $data[8,2] = 0.003; $data[8,3] = 0.199; $data[8,4] = 0.7; $data[8,5] = "Muhammad Iqbal Baqi"; $data[8,6] = "Benar"
$data[9,0] = "MIB_5.png"; $data[9,1] = 1.326; $data[9,2] = 0.004; $data[9,3] = 0.219; $data[9,4] = 0.667; $data[9,5] = "Muhammad Iqbal Baqi"; $data[9,6] = "Benar"
$data[10,0] = "AAH_1.png"; $data[10,1] = 0.737; $data[10,2] = 0.002; $data[10,3] = 0.24; $data[10,4] = 0.9; $data[10,5] = "Andrea Ayunove Hutami"; $data[10,6] = "Benar"
$data[11,0] = "AAH_2.png"; $data[11,1] = 0.951; $data[11,2] = 0.003; $data[11,3] = 0.6; $data[11,4] = 1; $data[11,5] = "Andrea Ayunove Hutami"; $data[11,6] = "Benar"
$data[12,0] = "AAH_3.png"; $data[12,1] = 0.715; $data[12,2] = 0.002; $data[12,3] = 0.244; $data[12,4] = 0.967; $data[12,5] = "Andrea Ayunove Hutami"; $data[12,6] = "Benar"
$data[13,0] = "TI_1.png"; $data[13,1] = 0.764; $data[13,2] = 0.003; $data[13,3] = 0.197; $data[13,4] = 0.533; $data[13,5] = "Muhammad Iqbal Baqi"; $data[13,6] = "Salah"
$data[14,0] = "TI_2.png"; $data[14,1] = 0.786; $data[14,2] = 0.003; $data[14,3] = 0.327; $data[14,4] = 0.533; $data[14,5] = "Toni Ismail"; $data[14,6] = "Benar"
$data[15,0] = "TI_3.png"; $data[15,1] = 0.554; $data[15,2] = 0.002; $data[15,3] = 0.554; $data[15,4] = 0.433; $data[15,5] = "Tidak Diketahui"; $data[15,6] = "Salah"
$data[16,0] = "TI_4.png"; $data[16,1] = 0.541; $data[16,2] = 0.002; $data[16,3] = 0.281; $data[16,4] = 0.6; $data[16,5] = "Toni Ismail"; $data[16,6] = "Benar"
$data[17,0] = "TI_5.png"; $data[17,1] = 0.913; $data[17,2] = 0.003; $data[17,3] = 0.356; $data[17,4] = 0.533; $data[17,5] = "Toni Ismail"; $data[17,6] = "Benar"
$data[18,0] = "RAS_1.png"; $data[18,1] = 0.477; $data[18,2] = 0.002; $data[18,3] = 0.371; $data[18,4] = 0.333; $data[18,5] = "Tidak Diketahui"; $data[18,6] = "Salah"
$data[19,0] = "RAS_2.png"; $data[19,1] = 0.855; $data[19,2] = 0.003; $data[19,3] = 0.269; $data[19,4] = 0.5; $data[19,5] = "Ridha Ayu Salsabila"; $data[19,6] = "Benar"
$data[20,0] = "RAS_3.png"; $data[20,1] = 0.358; $data[20,2] = 0.001; $data[20,3] = 0.343; $data[20,4] = 0.233; $data[20,5] = "Tidak Diketahui"; $data[20,6] = "Salah"
$data[21,0] = "RAS_4.png"; $data[21,1] = 1.082; $data[21,2] = 0.004; $data[21,3] = 0.097; $data[21,4] = 0.333; $data[21,5] = "Tidak Diketahui"; $data[21,6] = "Salah"
$data[22,0] = "RAS_5.png"; $data[22,1] = 1.008; $data[22,2] = 0.003; $data[22,3] = 0.345; $data[22,4] = 0.533; $data[22,5] = "Ridha Ayu Salsabila"; $data[22,6] = "Benar"
$data[23,0] = "RR_1.png"; $data[23,1] = 1.298; $data[23,2] = 0.004; $data[23,3] = 0.415; $data[23,4] = 0.667; $data[23,5] = "Rafiqo Rapitasari"; $data[23,6] = "Benar"
$data[24,0] = "RR_2.png"; $data[24,1] = 1.23; $data[24,2] = 0.004; $data[24,3] = 0.468; $data[24,4] = 0.5669999999999999; $data[24,5] = "Rafiqo Rapitasari"; $data[24,6] = "Benar"
$data[25,0] = "RR_3.png"; $data[25,1] = 0.928; $data[25,2] = 0.003; $data[25,3] = 0.096; $data[25,4] = 0.5669999999999999; $data[25,5] = "Rafiqo Rapitasari"; $data[25,6] = "Benar"
$data[26,0] = "RR_4.png"; $data[26,1] = 1.236; $data[26,2] = 0.004; $data[26,3] = 0.425; $data[26,4] = 0.633; $data[26,5] = "Rafiqo Rapitasari"; $data[26,6] = "Benar"
$data[27,0] = "RR_5.png"; $data[27,1] = 1.321; $data[27,2] = 0.004; $data[27,3] = 0.448; $data[27,4] = 0.733; $data[27,5] = "Rafiqo Rapitasari"; $data[27,6] = "Benar"
$data[28,0] = "AR_1.png"; $data[28,1] = 0.6870000000000001; $data[28,2] = 0.002; $data[28,3] = 0.317; $data[28,4] = 0.667; $data[28,5] = "Arizli Romadhon"; $data[28,6] = "Benar"
$data[29,0] = "GA_1.png"; $data[29,1] = 1.407; $data[29,2] = 0.005; $data[29,3] = 0.451; $data[29,4] = 0.9330000000000001; $data[29,5] = "Gege Ardiyansyah"; $data[29,6] = "Benar"
$data[30,0] = "GA_2.png"; $data[30,1] = 0.5; $data[30,2] = 0.002; $data[30,3] = 0.25; $data[30,4] = 0.7; $data[30,5] = "Gege Ardiyansyah"; $data[30,6] = "Benar"
$data[31,0] = "GA_3.png"; $data[31,1] = 0.5669999999999999; $data[31,2] = 0.002; $data[31,3] = 0.169; $data[31,4] = 0.6; $data[31,5] = "Gege Ardiyansyah"; $data[31,6] = "Benar"
$data[32,0] = "FY_1.png"; $data[32,1] = 0.91; $data[32,2] = 0.003; $data[32,3] = 0.232; $data[32,4] = 0.6; $data[32,5] = "Fanny Yusuf"; $data[32,6] = "Benar"
$data[33,0] = "FY_2.png"; $data[33,1] = 1.481; $data[33,2] = 0.005; $data[33,3] = 0.481; $data[33,4] = 0.3; $data[33,5] = "Tidak Diketahui"; $data[33,6] = "Salah"
$data[34,0] = "FY_3.png"; $data[34,1] = 1.376; $data[34,2] = 0.005; $data[34,3] = 0.297; $data[34,4] = 0.467; $data[34,5] = "Tidak Diketahui"; $data[34,6] = "Salah"
$data[35,0] = "FY_4.png"; $data[35,1] = 1.126; $data[35,2] = 0.004; $data[35,3] = 0.227; $data[35,4] = 0.367; $data[35,5] = "Tidak Diketahui"; $data[35,6] = "Salah"
$data[36,0] = "TO_1.png"; $data[36,1] = 0.633; $data[36,2] = 0.002; $data[36,3] = 0.357; $data[36,4] = 0.533; $data[36,5] = "Tiara Oktavian"; $data[36,6] = "Benar"
$data[37,0] = "TO_2.png"; $data[37,1] = 0.9429999999999999; $data[37,2] = 0.003; $data[37,3] = 0.368; $data[37,4] = 0.533; $data[37,5] = "Tiara Oktavian"; $data[37,6] = "Benar"
$data[38,0] = "TO_3.png"; $data[38,1] = 0.6860000000000001; $data[38,2] = 0.002; $data[38,3] = 0.391; $data[38,4] = 0.5; $data[38,5] = "Tiara Oktavian"; $data[38,6] = "Benar"
$data[39,0] = "TO_4.png"; $data[39,1] = 3.025; $data[39,2] = 0.01; $data[39,3] = 0.093; $data[39,4] = 0.6; $data[39,5] = "Tiara Oktavian"; $data[39,6] = "Benar"
$data[40,0] = "TO_5.png"; $data[40,1] = 2.877; $data[40,2] = 0.008999999999999999; $data[40,3] = 0.075; $data[40,4] = 0.633; $data[40,5] = "Tiara Oktavian"; $data[40,6] = "Benar"
$data[41,0] = "TD_1.png"; $data[41,1] = 2.154; $data[41,2] = 0.007; $data[41,3] = 0.029; $data[41,4] = 0.267; $data[41,5] = "Tidak Diketahui"; $data[41,6] = "Benar"
$data[42,0] = "TD_2.png"; $data[42,1] = 2.291; $data[42,2] = 0.008; $data[42,3] = 0.052; $data[42,4] = 0.3; $data[42,5] = "Tidak Diketahui"; $data[42,6] = "Benar"
$data[43,0] = "TD_3.png"; $data[43,1] = 0.919; $data[43,2] = 0.003; $data[43,3] = 0.172; $data[43,4] = 0.367; $data[43,5] = "Tidak Diketahui"; $data[43,6] = "Benar"
$data[44,0] = "TD_4.png"; $data[44,1] = 0.8100000000000001; $data[44,2] = 0.003; $data[44,3] = 0.097; $data[44,4] = 0.2; $data[44,5] = "Tidak Diketahui"; $data[44,6] = "Benar"

$ws.Range("A1:G45").Value = $data
